$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.331.52'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '1.621.76'
$ws.Range('E3').Value = '  +1.41%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '212.21'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').Value = '18.70'
$ws.Range('E10').Value = '  +2.27%  '
$ws.Range('D11').Value = '0.0814'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').Value = '1.847.59'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').Value = '1.625.33'
$ws.Range('E13').Value = '  +1.73%  '
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '26.341.60'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').Value = '62.43'
$ws.Range('E17').Value = '  +2.46%  '
$ws.Range('D18').Value = '0.0₃0724'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').Value = '202.10'
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('D22').Value = '9.29'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('E24').Value = '  -2.80%  '
$ws.Range('D25').Value = '144.30'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').Value = '15.19'
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').Value = '6.57'
$ws.Range('E29').Value = '  +0.47%  '
$ws.Range('E30').Value = '  +8.28%  '
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('E32').Value = '  +1.22%  '
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('E35').Value = '  +2.25%  '
$ws.Range('D36').Value = '1.160.42'
$ws.Range('E36').Value = '  +2.57%  '
$ws.Range('D37').Value = '0.0164'
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('D38').Value = '0.803'
$ws.Range('E38').Value = '  +0.60%  '
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').Value = '0.495'
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('D42').Value = '5.39'
$ws.Range('E42').Value = '  +3.97%  '
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').Value = '1.758.83'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('D45').Value = '92.31'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('D47').Value = '53.78'
$ws.Range('E47').Value = '  -0.83%  '
$ws.Range('E48').Value = '  +0.58%  '
$ws.Range('E49').Value = '  +0.59%  '
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').Value = '7.33'
$ws.Range('E51').Value = '  +1.87%  '
